$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 100001810
$ws.Cells.Item(29, 9).Value = 125000264
$ws.Cells.Item(29, 11).Value = 375000792
$ws.Cells.Item(29, 13).Value = -375000511
$ws.Cells.Item(33, 8).Value = 1140.45
$ws.Cells.Item(33, 9).Value = 159.07143
$ws.Cells.Item(33, 10).Value = 3430.3333
$ws.Cells.Item(33, 11).Value = 159.07143
$ws.Cells.Item(33, 12).Value = 3430.3333
$ws.Cells.Item(33, 13).Value = 69.92857000000001
$ws.Cells.Item(33, 14).Value = -3888.3333
$ws.Cells.Item(62, 8).Value = 64020
$ws.Cells.Item(62, 9).Value = 92920
$ws.Cells.Item(62, 11).Value = 92920
$ws.Cells.Item(62, 13).Value = -92296
$ws.Cells.Item(64, 8).Value = 7453.879
$ws.Cells.Item(64, 9).Value = 4097.9
$ws.Cells.Item(64, 10).Value = 8913
$ws.Cells.Item(64, 11).Value = 4097.9
$ws.Cells.Item(64, 12).Value = 8913
$ws.Cells.Item(64, 13).Value = -3849.9
$ws.Cells.Item(64, 14).Value = -9409
$ws.Cells.Item(65, 8).Value = 64020
$ws.Cells.Item(65, 9).Value = 92920
$ws.Cells.Item(65, 11).Value = 464600
$ws.Cells.Item(65, 13).Value = -461480
$ws.Cells.Item(67, 8).Value = 7453.879
$ws.Cells.Item(67, 9).Value = 4097.9
$ws.Cells.Item(67, 10).Value = 8913
$ws.Cells.Item(67, 11).Value = 4097.9
$ws.Cells.Item(67, 12).Value = 8913
$ws.Cells.Item(67, 13).Value = -3239.9
$ws.Cells.Item(67, 14).Value = -10629
$ws.Cells.Item(96, 8).Value = 770121.25
$ws.Cells.Item(96, 9).Value = 1111676.9
$ws.Cells.Item(96, 11).Value = 3335030.7
$ws.Cells.Item(96, 13).Value = -3333657.7
$ws.Cells.Item(113, 8).Value = 4198.0454
$ws.Cells.Item(113, 10).Value = 5398.6
$ws.Cells.Item(113, 12).Value = 5398.6
$ws.Cells.Item(113, 14).Value = -11906.6
$ws.Cells.Item(138, 8).Value = 1954.2452
$ws.Cells.Item(138, 10).Value = 2062.8206
$ws.Cells.Item(138, 12).Value = 6188.4618
$ws.Cells.Item(138, 14).Value = -16468.4618

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4999.857
$ws.Cells.Item(2, 9).Value = 4999
$ws.Cells.Item(2, 11).Value = 4999
$ws.Cells.Item(2, 13).Value = -4886
$ws.Cells.Item(63, 8).Value = 7250
$ws.Cells.Item(63, 9).Value = 4500
$ws.Cells.Item(63, 11).Value = 4500
$ws.Cells.Item(63, 13).Value = -3814
$ws.Cells.Item(66, 8).Value = 7250
$ws.Cells.Item(66, 9).Value = 4500
$ws.Cells.Item(66, 11).Value = 22500
$ws.Cells.Item(66, 13).Value = -19068
$ws.Cells.Item(74, 8).Value = 1930.7333
$ws.Cells.Item(74, 9).Value = 2088.182
$ws.Cells.Item(74, 11).Value = 2088.182
$ws.Cells.Item(74, 13).Value = -1214.182
$ws.Cells.Item(77, 8).Value = 1930.7333
$ws.Cells.Item(77, 9).Value = 2088.182
$ws.Cells.Item(77, 11).Value = 10440.91
$ws.Cells.Item(77, 13).Value = -6072.91
$ws.Cells.Item(88, 8).Value = 3800.6
$ws.Cells.Item(88, 9).Value = 3006
$ws.Cells.Item(88, 11).Value = 3006
$ws.Cells.Item(88, 13).Value = -2600
$ws.Cells.Item(91, 8).Value = 3800.6
$ws.Cells.Item(91, 9).Value = 3006
$ws.Cells.Item(91, 11).Value = 3006
$ws.Cells.Item(91, 13).Value = -1602
$ws.Cells.Item(116, 8).Value = 4999.857
$ws.Cells.Item(116, 9).Value = 4999
$ws.Cells.Item(116, 11).Value = 4999
$ws.Cells.Item(116, 13).Value = -2705
$ws.Cells.Item(122, 8).Value = 2195.5334
$ws.Cells.Item(122, 9).Value = 1975.2
$ws.Cells.Item(122, 10).Value = 2636.2
$ws.Cells.Item(122, 11).Value = 5925.6
$ws.Cells.Item(122, 12).Value = 7908.599999999999
$ws.Cells.Item(122, 13).Value = -3475.6
$ws.Cells.Item(122, 14).Value = -12808.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4999.857
$ws.Cells.Item(3, 9).Value = 4999
$ws.Cells.Item(3, 11).Value = 4999
$ws.Cells.Item(3, 13).Value = -4885
$ws.Cells.Item(20, 8).Value = 1606.5
$ws.Cells.Item(20, 9).Value = 1162
$ws.Cells.Item(20, 11).Value = 1162
$ws.Cells.Item(20, 13).Value = -915
$ws.Cells.Item(64, 8).Value = 658.4167
$ws.Cells.Item(64, 9).Value = 508.14285
$ws.Cells.Item(64, 10).Value = 868.8
$ws.Cells.Item(64, 11).Value = 508.14285
$ws.Cells.Item(64, 12).Value = 868.8
$ws.Cells.Item(64, 13).Value = -283.14285
$ws.Cells.Item(64, 14).Value = -1318.8
$ws.Cells.Item(67, 8).Value = 658.4167
$ws.Cells.Item(67, 9).Value = 508.14285
$ws.Cells.Item(67, 10).Value = 868.8
$ws.Cells.Item(67, 11).Value = 508.14285
$ws.Cells.Item(67, 12).Value = 868.8
$ws.Cells.Item(67, 13).Value = 271.85715
$ws.Cells.Item(67, 14).Value = -2428.8
$ws.Cells.Item(86, 8).Value = 5980
$ws.Cells.Item(86, 9).Value = 4133.3335
$ws.Cells.Item(86, 11).Value = 4133.3335
$ws.Cells.Item(86, 13).Value = -3010.3335
$ws.Cells.Item(89, 8).Value = 5980
$ws.Cells.Item(89, 9).Value = 4133.3335
$ws.Cells.Item(89, 11).Value = 20666.6675
$ws.Cells.Item(89, 13).Value = -15050.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 187.65218
$ws.Cells.Item(7, 9).Value = 28
$ws.Cells.Item(7, 10).Value = 436
$ws.Cells.Item(7, 11).Value = 28
$ws.Cells.Item(7, 12).Value = 436
$ws.Cells.Item(7, 13).Value = 85
$ws.Cells.Item(7, 14).Value = -662
$ws.Cells.Item(58, 8).Value = 1404.9412
$ws.Cells.Item(58, 9).Value = 1353.4546
$ws.Cells.Item(58, 11).Value = 1353.4546
$ws.Cells.Item(58, 13).Value = -1150.4546
$ws.Cells.Item(62, 8).Value = 7828.4287
$ws.Cells.Item(62, 9).Value = 933.3333
$ws.Cells.Item(62, 11).Value = 933.3333
$ws.Cells.Item(62, 13).Value = -309.3333
$ws.Cells.Item(65, 8).Value = 7828.4287
$ws.Cells.Item(65, 9).Value = 933.3333
$ws.Cells.Item(65, 11).Value = 4666.6665
$ws.Cells.Item(65, 13).Value = -1546.6665
$ws.Cells.Item(86, 8).Value = 6247.5
$ws.Cells.Item(86, 10).Value = 7000
$ws.Cells.Item(86, 12).Value = 7000
$ws.Cells.Item(86, 14).Value = -9246
$ws.Cells.Item(89, 8).Value = 6247.5
$ws.Cells.Item(89, 10).Value = 7000
$ws.Cells.Item(89, 12).Value = 35000
$ws.Cells.Item(89, 14).Value = -46232
$ws.Cells.Item(99, 8).Value = 2516.0908
$ws.Cells.Item(99, 9).Value = 2608.6667
$ws.Cells.Item(99, 10).Value = 2099.5
$ws.Cells.Item(99, 11).Value = 2608.6667
$ws.Cells.Item(99, 12).Value = 2099.5
$ws.Cells.Item(99, 13).Value = -1110.6667
$ws.Cells.Item(99, 14).Value = -5095.5
$ws.Cells.Item(122, 8).Value = 1380.25
$ws.Cells.Item(122, 9).Value = 1174.8334
$ws.Cells.Item(122, 10).Value = 1996.5
$ws.Cells.Item(122, 11).Value = 3524.5002
$ws.Cells.Item(122, 12).Value = 5989.5
$ws.Cells.Item(122, 13).Value = -1074.5002
$ws.Cells.Item(122, 14).Value = -10889.5
$ws.Cells.Item(126, 8).Value = 2516.0908
$ws.Cells.Item(126, 9).Value = 2608.6667
$ws.Cells.Item(126, 10).Value = 2099.5
$ws.Cells.Item(126, 11).Value = 7826.000100000001
$ws.Cells.Item(126, 12).Value = 6298.5
$ws.Cells.Item(126, 13).Value = -5356.000100000001
$ws.Cells.Item(126, 14).Value = -11238.5
$ws.Cells.Item(135, 8).Value = 200045.23
$ws.Cells.Item(135, 10).Value = 200045.23
$ws.Cells.Item(135, 12).Value = 200045.23
$ws.Cells.Item(135, 14).Value = -210185.23
$ws.Cells.Item(136, 8).Value = 1404.9412
$ws.Cells.Item(136, 9).Value = 1353.4546
$ws.Cells.Item(136, 11).Value = 4060.3638
$ws.Cells.Item(136, 13).Value = -1510.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 7339.8667
$ws.Cells.Item(23, 9).Value = 9880.333000000001
$ws.Cells.Item(23, 10).Value = 5646.222
$ws.Cells.Item(23, 11).Value = 29640.999
$ws.Cells.Item(23, 12).Value = 16938.666
$ws.Cells.Item(23, 13).Value = -29405.999
$ws.Cells.Item(23, 14).Value = -17408.666
$ws.Cells.Item(109, 8).Value = 97160.91
$ws.Cells.Item(109, 9).Value = 111991.11
$ws.Cells.Item(109, 10).Value = 30425
$ws.Cells.Item(109, 11).Value = 335973.33
$ws.Cells.Item(109, 12).Value = 91275
$ws.Cells.Item(109, 13).Value = -334933.33
$ws.Cells.Item(109, 14).Value = -93355
$ws.Cells.Item(122, 8).Value = 1687.8125
$ws.Cells.Item(122, 10).Value = 2295.3635
$ws.Cells.Item(122, 12).Value = 20658.2715
$ws.Cells.Item(122, 14).Value = -25558.2715
$ws.Cells.Item(131, 8).Value = 1672.6666
$ws.Cells.Item(131, 9).Value = 970.3333
$ws.Cells.Item(131, 10).Value = 2375
$ws.Cells.Item(131, 11).Value = 2910.9999
$ws.Cells.Item(131, 12).Value = 7125
$ws.Cells.Item(131, 13).Value = 2129.0001
$ws.Cells.Item(131, 14).Value = -17205

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4070.4062
$ws.Cells.Item(102, 9).Value = 2229.375
$ws.Cells.Item(102, 11).Value = 2229.375
$ws.Cells.Item(102, 13).Value = -607.375
$ws.Cells.Item(122, 8).Value = 41332.832
$ws.Cells.Item(122, 9).Value = 43599.4
$ws.Cells.Item(122, 11).Value = 130798.2
$ws.Cells.Item(122, 13).Value = -128348.2
$ws.Cells.Item(132, 8).Value = 2946.4443
$ws.Cells.Item(132, 9).Value = 2508.3333
$ws.Cells.Item(132, 10).Value = 3165.5
$ws.Cells.Item(132, 11).Value = 7524.999899999999
$ws.Cells.Item(132, 12).Value = 9496.5
$ws.Cells.Item(132, 13).Value = -4994.999899999999
$ws.Cells.Item(132, 14).Value = -14556.5
$ws.Cells.Item(134, 8).Value = 45789.8
$ws.Cells.Item(134, 10).Value = 45789.8
$ws.Cells.Item(134, 12).Value = 137369.4
$ws.Cells.Item(134, 14).Value = -142439.4
$ws.Cells.Item(136, 8).Value = 64500
$ws.Cells.Item(136, 10).Value = 64500
$ws.Cells.Item(136, 12).Value = 193500
$ws.Cells.Item(136, 14).Value = -198600

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 70000
$ws.Cells.Item(36, 10).Value = 70000
$ws.Cells.Item(36, 12).Value = 70000
$ws.Cells.Item(36, 14).Value = -71124
$ws.Cells.Item(122, 8).Value = 4585.4736
$ws.Cells.Item(122, 9).Value = 4864.3125
$ws.Cells.Item(122, 10).Value = 3098.3333
$ws.Cells.Item(122, 11).Value = 14592.9375
$ws.Cells.Item(122, 12).Value = 9294.999899999999
$ws.Cells.Item(122, 13).Value = -12142.9375
$ws.Cells.Item(122, 14).Value = -14194.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2962.4666
$ws.Cells.Item(81, 9).Value = 2328.4285
$ws.Cells.Item(81, 10).Value = 3517.25
$ws.Cells.Item(81, 11).Value = 4656.857
$ws.Cells.Item(81, 12).Value = 7034.5
$ws.Cells.Item(81, 13).Value = -3595.857
$ws.Cells.Item(81, 14).Value = -9156.5
$ws.Cells.Item(84, 8).Value = 2962.4666
$ws.Cells.Item(84, 9).Value = 2328.4285
$ws.Cells.Item(84, 10).Value = 3517.25
$ws.Cells.Item(84, 11).Value = 23284.285
$ws.Cells.Item(84, 12).Value = 35172.5
$ws.Cells.Item(84, 13).Value = -17980.285
$ws.Cells.Item(84, 14).Value = -45780.5
$ws.Cells.Item(122, 8).Value = 4262.696
$ws.Cells.Item(122, 9).Value = 2395.3845
$ws.Cells.Item(122, 10).Value = 6690.2
$ws.Cells.Item(122, 11).Value = 7186.1535
$ws.Cells.Item(122, 12).Value = 20070.6
$ws.Cells.Item(122, 13).Value = -4736.1535
$ws.Cells.Item(122, 14).Value = -24970.6
$ws.Cells.Item(136, 8).Value = 1799
$ws.Cells.Item(136, 9).Value = 448.5
$ws.Cells.Item(136, 11).Value = 1345.5
$ws.Cells.Item(136, 13).Value = 1204.5

